$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "522×9="
$t.Cell(1,2).Range.Text = "182×9="
$t.Cell(1,3).Range.Text = "251×3="
$t.Cell(1,4).Range.Text = "912×6="
$t.Cell(1,5).Range.Text = "363×4="

$t.Cell(5,1).Range.Text = "988×8="
$t.Cell(5,2).Range.Text = "613×2="
$t.Cell(5,3).Range.Text = "312×2="
$t.Cell(5,4).Range.Text = "257×7="
$t.Cell(5,5).Range.Text = "990×2="

$t.Cell(10,1).Range.Text = "291×8="
$t.Cell(10,2).Range.Text = "842×9="
$t.Cell(10,3).Range.Text = "826×2="
$t.Cell(10,4).Range.Text = "930×5="
$t.Cell(10,5).Range.Text = "395×7="

$t.Cell(15,1).Range.Text = "352×8="
$t.Cell(15,2).Range.Text = "291×6="
$t.Cell(15,3).Range.Text = "702×7="
$t.Cell(15,4).Range.Text = "657×2="
$t.Cell(15,5).Range.Text = "601×7="

$t.Cell(20,1).Range.Text = "568×8="
$t.Cell(20,2).Range.Text = "859×6="
$t.Cell(20,3).Range.Text = "835×9="
$t.Cell(20,4).Range.Text = "323×3="
$t.Cell(20,5).Range.Text = "735×7="
